# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit "Created functions to get season record": the sheet
# previously only had team/player stats (through column AC); now it also
# carries each player's team W-L-T record in columns AD:AF.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the style of the existing header row (bold,
# centered, thin border) by copying the format from the last header cell
# (AC1) and then writing the labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-51) gets the team's season record: 97 wins,
# 65 losses, 0 ties.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 0
}
